# Regenerate merged AHB files
# - Rename header suffixes "_old" -> "_FV2310" and "_new" -> "_FV2404"
# - Convert the data range into an Excel Table (ListObject) named "Table1"
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row values ---
# Columns A-J originally end in "_old" -> becomes "_FV2310"
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $current = [string]$cell.Value()
    $cell.Value() = ($current -replace "_old$", "_FV2310")
}

# Columns L-U originally end in "_new" -> becomes "_FV2404"
$newCols = @("L","M","N","O","P","Q","R","S","T","U")
foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $current = [string]$cell.Value()
    $cell.Value() = ($current -replace "_new$", "_FV2404")
}

# --- 2. Create an Excel Table (ListObject) over the full data range ---
$dataRange = $ws.Range("A1:U80")
$listObject = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$listObject.Name = "Table1"
$listObject.TableStyle() = ""

# --- 3. Freeze the header row ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes() = $true
